$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing 2022 column (J) to the new 2023 column (K)
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

# Fill in the new 2023 data
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1643.1
$ws.Range("K5").Value = 1158.7
$ws.Range("K6").Value = 1869.1
